# Chap05: End of copy pasting. Need actual writing.
# Update the cached "datetimeFigureOut" footer field text across the
# notes master, the slide master and all slide layouts (7/4/2017 -> 9/12/2017),
# and reposition/touch-up the stray "(c)" textbox on the single slide.

$p = $ppt.ActivePresentation

$oldDate = "7/4/2017"
$newDate = "9/12/2017"

function Update-DateField {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Notes master date placeholder.
Update-DateField $p.NotesMaster.Shapes

# Slide master date placeholder.
Update-DateField $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateField $layouts.Item($j).Shapes
}

# Slide 1: nudge the "(c)" caption textbox up
# (the copy-pasted placement left behind by the previous edit).
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "ZoneTexte 7") {
        $sh.Top = 344.4550393700787
    }
}
